$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4427559971809387
$ws.Range("B1").Value = 0.5897301435470581
$ws.Range("C1").Value = 0.9554622173309326
$ws.Range("D1").Value = 5.439185619354248
$ws.Range("E1").Value = 3.173203706741333
